# Regenerate the localization-status report for the renamed source file:
#   ee778694-5ccc-4a6d-9240-5f72d64a0014.md -> ce8acbdd-cf47-4107-b5fc-e315123a475f.md
# A new handoff xliff was generated for both locales and the previously
# recorded handback info is cleared out (report is "Ready for handoff").

$wb = $excel.ActiveWorkbook

$oldGuid = "ee778694-5ccc-4a6d-9240-5f72d64a0014"
$newGuid = "ce8acbdd-cf47-4107-b5fc-e315123a475f"
$oldHash = "56a6c5df332dc760f00010a0507bc54ab8f66aa3"
$newHash = "78eede814804593c668d22ec84653ebfdd61e536"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = "$newGuid.md"
$overview.Range("B2").Value = "e2e\$newGuid.md"
$overview.Range("G2").Value = "2016-08-28 20:58:31"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = "$newGuid.md"
$zhcn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-28 20:58:27"
$zhcn.Range("I2").Value = ""
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = "$newGuid.md"
$dede.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$dede.Range("H2").Value = "2016-08-28 20:58:31"
$dede.Range("I2").Value = ""
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = "0001-01-01 00:00:00"
